# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right before the existing "2022-Q3"
#    sheet (i.e. right after "总计"), populated with the new quarter's
#    fund-holding data. We duplicate the "2022-Q3" sheet so the new sheet
#    inherits the same layout/styles (header row + index-column formatting),
#    then overwrite the cell contents.
# 2. Insert a new top data row in the "总计" (summary) sheet for 2022-Q4,
#    pushing the existing 2022-Q3 / 2021-Q2 / 2021-Q1 rows down by one.

$wb = $excel.ActiveWorkbook

# Helper: write a value as plain TEXT (no leading-apostrophe / quote-prefix
# style baked in) even when it looks like a number - mirrors how the source
# file stores fund codes / percentages as inline strings. We do this by
# writing a formula that evaluates to the literal text, then Copy +
# PasteSpecial(values) to collapse it back down to a static value, which
# leaves the cell's style untouched (no quotePrefix xf gets attached).
function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = "=""$escaped"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q4" sheet before "2022-Q3" (clone it so the
# header/index-column styling matches the sibling quarter sheets).
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "2022-Q4"

# Clear out the old (copied) data rows beyond the header so we start clean.
$newSheet.Rows.Item(2).Resize(10).ClearContents()

$newSheet.Cells.Item(2,1).Value = 0
Set-TextValue $newSheet.Cells.Item(2,2) "013442"
$newSheet.Cells.Item(2,3).Value = "建信中证1000指数增强E"
Set-TextValue $newSheet.Cells.Item(2,4) "9.52"
Set-TextValue $newSheet.Cells.Item(2,5) "86.80"
Set-TextValue $newSheet.Cells.Item(2,6) "1.38"
Set-TextValue $newSheet.Cells.Item(2,7) "0.1314"
$newSheet.Cells.Item(2,8).Value = 4

$newSheet.Cells.Item(3,1).Value = 1
Set-TextValue $newSheet.Cells.Item(3,2) "006165"
$newSheet.Cells.Item(3,3).Value = "建信中证1000指数增强A"
Set-TextValue $newSheet.Cells.Item(3,4) "7.20"
Set-TextValue $newSheet.Cells.Item(3,5) "86.80"
Set-TextValue $newSheet.Cells.Item(3,6) "1.38"
Set-TextValue $newSheet.Cells.Item(3,7) "0.0994"
$newSheet.Cells.Item(3,8).Value = 4

$newSheet.Cells.Item(4,1).Value = 2
Set-TextValue $newSheet.Cells.Item(4,2) "006166"
$newSheet.Cells.Item(4,3).Value = "建信中证1000指数增强C"
Set-TextValue $newSheet.Cells.Item(4,4) "2.21"
Set-TextValue $newSheet.Cells.Item(4,5) "86.80"
Set-TextValue $newSheet.Cells.Item(4,6) "1.38"
Set-TextValue $newSheet.Cells.Item(4,7) "0.0305"
$newSheet.Cells.Item(4,8).Value = 4

# The index column (A) needs the same centred/bordered style as the header
# row and as the sibling sheets' index columns; pick it up from A1's row-2
# neighbour on the source sheet via copy/paste-format.
$idxStyleSrc = $q3Sheet.Cells.Item(2,1)
$idxStyleSrc.Copy()
$newSheet.Range($newSheet.Cells.Item(2,1), $newSheet.Cells.Item(4,1)).PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: shift the "总计" summary rows down and insert the 2022-Q4 row
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 5 (2021-Q1, previously row 4) needs the same index-column style as
# the existing rows above it before we fill its value in.
$summary.Cells.Item(4,1).Copy()
$summary.Cells.Item(5,1).PasteSpecial(-4122)

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q1"
$summary.Cells.Item(5,3).Value = 3
$summary.Cells.Item(5,4).Value = 0.06

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2021-Q2"
$summary.Cells.Item(4,3).Value = 3
$summary.Cells.Item(4,4).Value = 0.09

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q3"
$summary.Cells.Item(3,3).Value = 2
$summary.Cells.Item(3,4).Value = 0.03

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 0.26
